# Add a new trade record (row 4) to the sheet, mirroring the layout/format
# of the existing row 3 (e.g. the date column G keeps its date number format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 3's formatting down into row 4 first so new cells (in particular
# the date cell G4) pick up the same styles (e.g. style index for dates).
$ws.Range("A3:H3").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)  # xlPasteFormats

# Now populate the new row's values.
$ws.Range("A4").Value = 10035.5
$ws.Range("B4").Value = 9945
$ws.Range("C4").Value = 107.96
$ws.Range("D4").Value = 108.94
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.91
$ws.Range("G4").Value = 42609.505567129629
$ws.Range("H4").Value = $true
